# Apply the cryptos-list update (GitHub Actions data refresh) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.091.33"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.652.11"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.35"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5273"
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2599"
$ws.Range("E8").Value = "  -1.52%  "
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("E10").Value = "  -2.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07785"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").Value = "1.632.34"
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("D14").Value = "1.879.28"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5496"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "0.0₅8214"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.55"
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("D18").Value = "26.089.90"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.579"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.82"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.035"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.03"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1239"
$ws.Range("E26").Value = "  +1.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.231"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.431"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05847"
$ws.Range("E30").Value = "  -1.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.273"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.548"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.268"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.585"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9453"
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.409"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5740"
$ws.Range("E38").Value = "  +1.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01609"
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8432"
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.741"
$ws.Range("E41").Value = "  -4.88%  "
$ws.Range("E42").Value = "  +2.88%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "1.030.27"
$ws.Range("E44").Value = "  +1.49%  "
$ws.Range("D45").Value = "1.794.96"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.06"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.004"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("E48").Value = "  +2.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05142"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.817"
$ws.Range("E50").Value = "  -3.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.456"
$ws.Range("E51").Value = "  +0.41%  "

Write-Host "Updated cryptos list"
